$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")
$ws.Rows.Item(16).Delete()
$ws.Rows.Item(13).Delete()
